# nexial-core "#system" sheet (sheet1) drives the dropdown lists used on the
# io_showcase sheet via named ranges. The "web" named range (column U,
# U2:U111) lists every available `web` command; this change adds the new
# dragAndDrop(fromLocator,toLocator) command to that alphabetically sorted
# list and grows the named range by one row (U2:U112).
#
# "dragAndDrop" sorts between "doubleClickByLabelAndWait(label,waitMs)" and
# "editLocalStorage(key,value)", i.e. it belongs at row 59 - so everything
# currently at U59:U111 shifts down one row to U60:U112, and the new command
# is written into the now-empty U59.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# Shift the existing tail of the "web" list down by one row, then insert
# the new command at the vacated slot. Only column U is touched - other
# columns (other command lists such as "desktop" in column E) must be left
# completely alone, so this uses a plain column range copy instead of a
# row/entire-row insert.
$tail = $ws.Range("U59:U111").Value2
$ws.Range("U60:U112").Value = $tail
$ws.Range("U59").Value = "dragAndDrop(fromLocator,toLocator)"

# Grow the "web" defined name so it covers the newly added row.
$wb.Names.Item("web").RefersTo = "='#system'!`$U`$2:`$U`$112"
